# Pizza price workbook clean-up edit
# - Corrects the D9 price from 453 to 433
# - Moves the sheet's active selection from A10 to F10
# - Updates the theme's "Background 1" / window color from white (FFFFFF) to A8A8A8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the price value in D9 (was 453, should be 433)
$ws.Range("D9").Value = 433

# Update the saved cell selection/active cell to F10
$ws.Range("F10").Select() | Out-Null

# Update theme color: Background 1 (lt1) from FFFFFF to A8A8A8
$themeColors = $wb.Theme.ThemeColorScheme
$background1 = $themeColors.Colors(2)
$background1.RGB = 0xA8A8A8
